$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1364
$ws1.Range("F3").Value = 1373
$ws1.Range("F5").Value = 116
$ws1.Range("F6").Value = 70
$ws1.Range("F7").Value = 678
$ws1.Range("F8").Value = 117
$ws1.Range("F9").Value = 52
$ws1.Range("F11").Value = 2467
$ws1.Range("F12").Value = 1607
$ws1.Range("F13").Value = 1501
$ws1.Range("F14").Value = 312
$ws1.Range("F15").Value = 247
$ws1.Range("F16").Value = 617
$ws1.Range("F17").Value = 789
$ws1.Range("F18").Value = 83
$ws1.Range("F22").Value = 29
$ws1.Range("F24").Value = 5065
$ws1.Range("F25").Value = 220
$ws1.Range("F26").Value = 539
$ws1.Range("F27").Value = 81
$ws1.Range("F29").Value = 138
$ws1.Range("F30").Value = 227
$ws1.Range("F31").Value = 224
$ws1.Range("F32").Value = 32
$ws1.Range("F33").Value = 1041
$ws1.Range("F34").Value = 739
$ws1.Range("F36").Value = 53
$ws1.Range("F39").Value = 1075
$ws1.Range("F42").Value = 174
$ws1.Range("F44").Value = 44
$ws1.Range("G14").Value = 218

# Sheet 2: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 7

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1364
$ws4.Range("F5").Value = 1373
$ws4.Range("F9").Value = 116
$ws4.Range("F10").Value = 70
$ws4.Range("F11").Value = 678
$ws4.Range("F12").Value = 117
$ws4.Range("F13").Value = 7
$ws4.Range("F17").Value = 2467
$ws4.Range("F18").Value = 1607
$ws4.Range("F19").Value = 1501
$ws4.Range("F20").Value = 312
$ws4.Range("F21").Value = 247
$ws4.Range("F22").Value = 617
$ws4.Range("F24").Value = 789
$ws4.Range("F25").Value = 83
$ws4.Range("F28").Value = 29
$ws4.Range("F29").Value = 5065
$ws4.Range("F30").Value = 220
$ws4.Range("F31").Value = 539
$ws4.Range("F32").Value = 81
$ws4.Range("F34").Value = 138
$ws4.Range("F35").Value = 227
$ws4.Range("F36").Value = 224
$ws4.Range("F37").Value = 32
$ws4.Range("F38").Value = 1041
$ws4.Range("F39").Value = 739
$ws4.Range("F40").Value = 53
$ws4.Range("F42").Value = 1075
$ws4.Range("F44").Value = 174
$ws4.Range("F46").Value = 44
$ws4.Range("G20").Value = 218
